# rggen/rggen#14 -- register_8 (0x20) gains a third bit field.
#
# Before: register_8 had 2 bit fields (bit_field_0 @ rows 42, bit_field_1 @ row 43)
# After:  register_8 has 3 bit fields (bit_field_0 @ row 42, bit_field_1 @ row 43 [values
#         updated], bit_field_2 @ new row 44). Everything below register_8 shifts down
#         by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) bit_field_0's bit-range notation changes (register now has 4 fields of width 4 each
#    instead of 2 fields of width 8, so the "msb:width:num-fields:reg-width" annotation
#    changes from 0:8:4:16 to 0:4:4:16).
$ws.Range("G42").Value() = "0:4:4:16"

# 2) Insert a new row before the old row 44 (register_9 header), so the new bit_field_2
#    row lands at row 44 and everything from the old row 44 onward shifts down by one.
$ws.Rows.Item(44).Insert()

# 3) The new row 44 should look like the rest of the register_8 bit-field rows (borders
#    etc.) -- grab formatting from row 43 (currently still the last bit-field row of
#    register_8, so it already carries the "last row of block" bottom border).
$ws.Range("B43:J43").Copy()
$ws.Range("B44:J44").PasteSpecial(-4122)

# 4) Row 43 is no longer the last bit-field row of register_8 (row 44 now is), so its
#    B:E block-border formatting needs to switch from "last row" to "middle row" style.
#    Row 41 is a plain middle row with that exact formatting -- copy it.
$ws.Range("B41:E41").Copy()
$ws.Range("B43:E43").PasteSpecial(-4122)

# 5) Fill in the content for the (updated) bit_field_1 row and the (new) bit_field_2 row.
$ws.Range("F43").Value() = "bit_field_1"
$ws.Range("G43").Value() = "4:4:4:16"
$ws.Range("H43").Value() = "rw"
$ws.Range("I43").Value() = "default: 0"

$ws.Range("F44").Value() = "bit_field_2"
$ws.Range("G44").Value() = "8:4:4:16"
$ws.Range("H44").Value() = "rw"
$ws.Range("I44").Value() = "0, 1, 2, 3"
